$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the Specifications intro
#    paragraph (it will be re-added later after the Exclusions text).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Fill in the Inclusions placeholder.
# ------------------------------------------------------------------
$inclPara = $d.Paragraphs.Item(11)
$inclPara.Range.Find.Execute("TODO", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Messages listed in both RFCs", 2)

# ------------------------------------------------------------------
# 3. Fill in the Exclusions placeholder, then re-create the "_GoBack"
#    bookmark right after the new text (collapsed, not spanning it).
#    A trailing placeholder character is used to dodge an end-of-
#    paragraph edge case in the bookmark placement logic, and is
#    removed immediately afterwards.
# ------------------------------------------------------------------
$exclPara = $d.Paragraphs.Item(14)
$exclPara.Range.Find.Execute("TODO", $true, $false, $false, $false, $false, `
    $true, 1, $false, "No exclusionsZ", 2)

$exclPara = $d.Paragraphs.Item(14)
$bmPos = $exclPara.Range.Start + "No exclusions".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()

# ------------------------------------------------------------------
# 4. Merge the three runs of the "Test Environment" paragraph back
#    into a single run, dropping the gramStart/gramEnd proofErr
#    markers that bracketed "supported". Replacing a short run of
#    text that straddles the run boundaries (but doesn't touch the
#    apostrophe, to avoid it being smart-quoted) is enough to make
#    Word re-flow the paragraph into one run.
# ------------------------------------------------------------------
$testEnvRange = $d.Paragraphs.Item(98).Range
$testEnvRange.Find.Execute("is supported, Huawei", $true, $false, $false, `
    $false, $false, $true, 1, $false, "is supported, Huawei", 2)
